# Adds a new forecast-date column (Q) and a new observed-date row (29) to
# both the "cases" and "deaths" tables, and backfills the newly-observed
# value for 2020-04-18 (row 15) in the "Observed" column (B).
#
# Table shape reminder: column A holds the "observed" date for the row;
# row 1 (columns C..) holds the date the forecast was produced on; column
# B is the actually-observed count for that date. Each row's forecasts
# walk diagonally rightwards across C..P (one new forecast column per
# day), so adding a new forecast date adds one new column (Q) and one new
# diagonal value per existing row (from the row whose forecast-of-record
# had run off the right edge of the table onward), plus a brand new row
# for the newly observed date itself.

$wb = $excel.ActiveWorkbook

# New forecast date header used for column Q (continues the header-row
# sequence ...,"2020-04-16","2020-04-17","2020-04-18").
$newForecastDate = "2020-04-18"

# New observed date used for the new row 29 (continues the A-column
# sequence ...,"2020-04-30","2020-05-01","2020-05-02").
$newObservedDate = "2020-05-02"

# Per-sheet data: the backfilled B15 ("Observed" on 2020-04-18") value,
# and the new diagonal Q values for rows 16..29 (2020-04-19 .. 2020-05-02).
$sheetData = [ordered]@{
    "cases"  = @{
        B15 = 36599
        Q   = [ordered]@{
            16 = 38690
            17 = 41351
            18 = 44612
            19 = 47099
            20 = 49109
            21 = 51377
            22 = 54697
            23 = 56663
            24 = 60947
            25 = 63232
            26 = 66951
            27 = 70303
            28 = 73768
            29 = 76688
        }
    }
    "deaths" = @{
        B15 = 2347
        Q   = [ordered]@{
            16 = 2513
            17 = 2714
            18 = 2964
            19 = 3157
            20 = 3314
            21 = 3492
            22 = 3755
            23 = 3912
            24 = 4257
            25 = 4443
            26 = 4747
            27 = 5024
            28 = 5312
            29 = 5557
        }
    }
}

function Set-TextValue($range, $text) {
    # Writing a date-shaped literal like "2020-04-18" through .Value lets
    # Excel's normal autoconvert turn it into a real date serial (plus a
    # date number-format). The source table stores these as plain shared
    # strings instead, so force text entry via a "@" number format and
    # then drop the formatting again so the cell is left as plain text
    # with no explicit style, matching the rest of the date-label column.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

foreach ($name in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $data = $sheetData[$name]

    # New column Q: header (row 1) is the next forecast date.
    Set-TextValue $ws.Range("Q1") $newForecastDate

    # Backfill the newly observed count for 2020-04-18 (row 15, column B).
    $ws.Range("B15").Value = $data.B15

    # New diagonal forecast values in column Q for rows 16..29.
    foreach ($r in $data.Q.Keys) {
        $ws.Cells.Item($r, 17).Value = $data.Q[$r]
    }

    # New row 29: observed date in column A, diagonal forecast in Q (set
    # above via the Q loop already covers row 29's value).
    Set-TextValue $ws.Range("A29") $newObservedDate
}
